$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row stays the same text ("day" / "name_duty"); shared-string index
# shifts are handled automatically when the workbook is re-saved.
$ws.Range("A1").Value = "day"
$ws.Range("B1").Value = "name_duty"

# New duty roster values for column B (rows 2-32)
$ws.Range("B2").Value = "三神佳誠"
$ws.Range("B3").Value = "氏家琉貴"
$ws.Range("B4").Value = "羽賀尚生"
$ws.Range("B5").Value = "島田実"
$ws.Range("B6").Value = "足立耕平"
$ws.Range("B7").Value = "遠藤隼人"
$ws.Range("B8").Value = "Ethan Virtudazo"
$ws.Range("B9").Value = "富澤天音"
$ws.Range("B10").Value = "神山修造"
$ws.Range("B11").Value = "川田涼介"
$ws.Range("B12").Value = ""
$ws.Range("B13").Value = "豊島亮"
$ws.Range("B14").Value = "兒島大志郎"
$ws.Range("B15").Value = "高野怜央"
$ws.Range("B16").Value = "山口玲"
$ws.Range("B17").Value = "日高泰聖"
$ws.Range("B18").Value = "志塚惇希"
$ws.Range("B19").Value = "山口洸翔"
$ws.Range("B20").Value = "白岩詩佑介"
$ws.Range("B21").Value = "石井海成"
$ws.Range("B22").Value = "Nicholas Tristan Aryasatyo"
$ws.Range("B23").Value = "小溝賢"
$ws.Range("B24").Value = "小野文哉"
$ws.Range("B25").Value = "渡部魁"
$ws.Range("B26").Value = "崎谷航平, Jun Seomun"
$ws.Range("B27").Value = "三神佳誠"
$ws.Range("B28").Value = "氏家琉貴"
$ws.Range("B29").Value = "羽賀尚生"
$ws.Range("B30").Value = "島田実"
$ws.Range("B31").Value = "足立耕平"
$ws.Range("B32").Value = "遠藤隼人"

# "Ethan Virtudazo" keeps the small Roboto font used for Latin names
$ws.Range("B8").Font.Name = "Roboto"
$ws.Range("B8").Font.Size = 10

# The row that previously held "Ethan Virtudazo" (and its Roboto font)
# now holds a regular Japanese name again, so restore the plain Arial font
$ws.Range("B17").Font.Name = "Arial"
$ws.Range("B17").Font.Size = 10

# Update the saved selection/active cell
$ws.Range("C16").Select()
